$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("02 Oct (1cm Interval)")
$ws2.Range("B11").Value = 455.6
